$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.003254057381738562
$ws.Range("D2").Value = 0.03349782192563566
$ws.Range("E2").Value = 0.4333214853876086
$ws.Range("F2").Value = 1.479172793580574
$ws.Range("G2").Value = 1.560190433613002
$ws.Range("H2").Value = 1.022540459300615
$ws.Range("I2").Value = 1.1362179686612
$ws.Range("C3").Value = 0.002821308720768911
$ws.Range("D3").Value = 0.02934618752100704
$ws.Range("E3").Value = 0.3770546871430156
$ws.Range("F3").Value = 1.338957110356432
$ws.Range("G3").Value = 1.390408734690794
$ws.Range("H3").Value = 0.9498434618274416
$ws.Range("I3").Value = 1.027328598235428
$ws.Range("C4").Value = 0.002556772847700728
$ws.Range("D4").Value = 0.02681457152093714
$ws.Range("E4").Value = 0.3427093422305205
$ws.Range("F4").Value = 1.25385654906384
$ws.Range("G4").Value = 1.287088664196489
$ws.Range("H4").Value = 0.9059721192086272
$ws.Range("I4").Value = 0.9612294625744511
$ws.Range("C5").Value = 0.002449222314204746
$ws.Range("D5").Value = 0.02578698926784995
$ws.Range("E5").Value = 0.3287593870395824
$ws.Range("F5").Value = 1.219418417453852
$ws.Range("G5").Value = 1.245207559488222
$ws.Range("H5").Value = 0.8882818619924819
$ws.Range("I5").Value = 0.9344778148814612
$ws.Range("C6").Value = 0.002431377428983694
$ws.Range("D6").Value = 0.02561659642820757
$ws.Range("E6").Value = 0.32644564704799
$ws.Range("F6").Value = 1.213714304211692
$ws.Range("G6").Value = 1.23826638392859
$ws.Range("H6").Value = 0.885355608407167
$ws.Range("I6").Value = 0.930046662266534
$ws.Range("C7").Value = 0.002555321430154578
$ws.Range("D7").Value = 0.02680069713186128
$ws.Range("E7").Value = 0.3425210282276652
$ws.Range("F7").Value = 1.25339113956133
$ws.Range("G7").Value = 1.286522952935513
$ws.Range("H7").Value = 0.9057327884428616
$ws.Range("I7").Value = 0.9608679430434961
$ws.Range("C8").Value = 0.003104574512214242
$ws.Range("D8").Value = 0.03206250382659448
$ws.Range("E8").Value = 0.4138756201795957
$ws.Range("F8").Value = 1.430615228876633
$ws.Range("G8").Value = 1.50145106635506
$ws.Range("H8").Value = 0.997312915528255
$ws.Range("I8").Value = 1.098511222843172
$ws.Range("C9").Value = 0.004193089834402031
$ws.Range("D9").Value = 0.04253537497446302
$ws.Range("E9").Value = 0.5556417919495686
$ws.Range("F9").Value = 1.786433726092952
$ws.Range("G9").Value = 1.930770471932931
$ws.Range("H9").Value = 1.183194529484183
$ws.Range("I9").Value = 1.374773752377479
$ws.Range("C10").Value = 0.005002855156178043
$ws.Range("D10").Value = 0.05034663962426578
$ws.Range("E10").Value = 0.6612548373531837
$ws.Range("F10").Value = 2.053503134825348
$ws.Range("G10").Value = 2.25169995454138
$ws.Range("H10").Value = 1.323928293468384
$ws.Range("I10").Value = 1.582080857590341
$ws.Range("C11").Value = 0.005374107612826862
$ws.Range("D11").Value = 0.05393065389280594
$ws.Range("E11").Value = 0.7096928564080969
$ws.Range("F11").Value = 2.176361795961952
$ws.Range("G11").Value = 2.399058331659944
$ws.Range("H11").Value = 1.388930860006212
$ws.Range("I11").Value = 1.677437594168651
$ws.Range("C12").Value = 0.005515163530720457
$ws.Range("D12").Value = 0.05529264728673411
$ws.Range("E12").Value = 0.7280979333847029
$ws.Range("F12").Value = 2.223092649497147
$ws.Range("G12").Value = 2.455068630089329
$ws.Range("H12").Value = 1.413692841925524
$ws.Range("I12").Value = 1.713706417683738
$ws.Range("C13").Value = 0.005484762777026475
$ws.Range("D13").Value = 0.05499909744283116
$ws.Range("E13").Value = 0.7241311873571874
$ws.Range("F13").Value = 2.213018952958919
$ws.Range("G13").Value = 2.442996318239921
$ws.Range("H13").Value = 1.408353279212236
$ws.Range("I13").Value = 1.705888062043613
$ws.Range("C14").Value = 0.005385702609842724
$ws.Range("D14").Value = 0.05404260727604537
$ws.Range("E14").Value = 0.7112057626939219
$ws.Range("F14").Value = 2.180202163862305
$ws.Range("G14").Value = 2.403662073796909
$ws.Range("H14").Value = 1.390965065424609
$ws.Range("I14").Value = 1.680418212570913
$ws.Range("C15").Value = 0.005325088373179909
$ws.Range("D15").Value = 0.05345736725728045
$ws.Range("E15").Value = 0.7032969061931595
$ws.Range("F15").Value = 2.160128194102469
$ws.Range("G15").Value = 2.379596291095424
$ws.Range("H15").Value = 1.38033358776039
$ws.Range("I15").Value = 1.664838182651778
$ws.Range("C16").Value = 0.004978655728585579
$ws.Range("D16").Value = 0.05011306619309153
$ws.Range("E16").Value = 0.6580977365154581
$ws.Range("F16").Value = 2.045502423609747
$ws.Range("G16").Value = 2.242098281615256
$ws.Range("H16").Value = 1.319700468231076
$ws.Range("I16").Value = 1.575870917717282
$ws.Range("C17").Value = 0.004766908780645451
$ws.Range("D17").Value = 0.04806957234971776
$ws.Range("E17").Value = 0.6304746509790249
$ws.Range("F17").Value = 1.975540707466763
$ws.Range("G17").Value = 2.158106315897214
$ws.Range("H17").Value = 1.282759550197568
$ws.Range("I17").Value = 1.521567415056438
$ws.Range("C18").Value = 0.004645385760881027
$ws.Range("D18").Value = 0.04689706723206655
$ws.Range("E18").Value = 0.6146233109359969
$ws.Range("F18").Value = 1.935428655982946
$ws.Range("G18").Value = 2.109924162734671
$ws.Range("H18").Value = 1.261604197832753
$ws.Range("I18").Value = 1.490431971508826
$ws.Range("C19").Value = 0.004604284556375404
$ws.Range("D19").Value = 0.04650055608500736
$ws.Range("E19").Value = 0.6092624483873266
$ws.Range("F19").Value = 1.921869092488407
$ws.Range("G19").Value = 2.093632089949608
$ws.Range("H19").Value = 1.254457003917366
$ws.Range("I19").Value = 1.479906719779933
$ws.Range("C20").Value = 0.004789421446250941
$ws.Range("D20").Value = 0.04828680705054467
$ws.Range("E20").Value = 0.6334113316135301
$ws.Range("F20").Value = 1.98297492231211
$ws.Range("G20").Value = 2.167034087632203
$ws.Range("H20").Value = 1.286682402595602
$ws.Range("I20").Value = 1.527337864055909
$ws.Range("C21").Value = 0.005414785755259999
$ws.Range("D21").Value = 0.05432341797123286
$ws.Range("E21").Value = 0.7150005244285893
$ws.Range("F21").Value = 2.189835550321789
$ws.Range("G21").Value = 2.415209724981651
$ws.Range("H21").Value = 1.396068371685772
$ws.Range("I21").Value = 1.687894936110695
$ws.Range("C22").Value = 0.005826268009457181
$ws.Range("D22").Value = 0.05829689329300436
$ws.Range("E22").Value = 0.7686918527812168
$ws.Range("F22").Value = 2.326242151678684
$ws.Range("G22").Value = 2.578630384526036
$ws.Range("H22").Value = 1.468417554920734
$ws.Range("I22").Value = 1.793760671793962
$ws.Range("C23").Value = 0.005606379835924713
$ws.Range("D23").Value = 0.056173459059778
$ws.Range("E23").Value = 0.740000080994875
$ws.Range("F23").Value = 2.253325073928437
$ws.Range("G23").Value = 2.49129354428328
$ws.Range("H23").Value = 1.429722897787428
$ws.Range("I23").Value = 1.737170097874895
$ws.Range("C24").Value = 0.004779242810620588
$ws.Range("D24").Value = 0.04818858802501325
$ws.Range("E24").Value = 0.6320835667388707
$ws.Range("F24").Value = 1.979613572332795
$ws.Range("G24").Value = 2.162997512968275
$ws.Range("H24").Value = 1.284908624568857
$ws.Range("I24").Value = 1.524728781897011
$ws.Range("C25").Value = 0.003897082456209233
$ws.Range("D25").Value = 0.03968333960929726
$ws.Range("E25").Value = 0.5170591309376391
$ws.Range("F25").Value = 1.68922753228091
$ws.Range("G25").Value = 1.813714311031788
$ws.Range("H25").Value = 1.13220129104684
$ws.Range("I25").Value = 1.299310391377276
